$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Regroup deplucated proprietaire in <annex 1>": the "HAYLALA ONE" proprietaire
# had a contract (949/DR) already listed further down; this inserts two more
# of their contracts right above it so all of HAYLALA ONE's rows sit together
# at the top of the table (rows 2-4), pushing everything else down by 2 rows.
$ws.Rows.Item(2).Insert() | Out-Null
$ws.Rows.Item(2).Insert() | Out-Null

# New row 2: HAYLALA ONE - contract 049/DR 2 - annuelle
$ws.Cells.Item(2, 1).Value = "HAYLALA ONE"
$ws.Cells.Item(2, 2).Value = "BG12456"
$ws.Cells.Item(2, 3).Value = "'111111111111111111111111"
$ws.Cells.Item(2, 4).Value = "'1"
$ws.Cells.Item(2, 5).Value = "'11"
$ws.Cells.Item(2, 6).Value = "Direction régionale"
$ws.Cells.Item(2, 7).Value = "049/DR 2"
$ws.Cells.Item(2, 8).Value = "annuelle"
$ws.Cells.Item(2, 9).Value = 200000
$ws.Cells.Item(2, 10).Value = 10000
$ws.Cells.Item(2, 11).Value = 190000

# New row 3: HAYLALA ONE - contract 094/DR 1 - mensuelle
$ws.Cells.Item(3, 1).Value = "HAYLALA ONE"
$ws.Cells.Item(3, 2).Value = "BG12456"
$ws.Cells.Item(3, 3).Value = "'111111111111111111111111"
$ws.Cells.Item(3, 4).Value = "'1"
$ws.Cells.Item(3, 5).Value = "'11"
$ws.Cells.Item(3, 6).Value = "Direction régionale"
$ws.Cells.Item(3, 7).Value = "094/DR 1"
$ws.Cells.Item(3, 8).Value = "mensuelle"
$ws.Cells.Item(3, 9).Value = 20000
$ws.Cells.Item(3, 10).Value = 1500
$ws.Cells.Item(3, 11).Value = 18500

# Update totals row (now row 8, was row 6) with new sums that include the two
# newly added rows.
$ws.Cells.Item(8, 9).Value = 280000.13
$ws.Cells.Item(8, 10).Value = 15500.01
$ws.Cells.Item(8, 11).Value = 264500.12
